$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 2) ----
$ws.Range("B2").Value = "Tâches"
$ws.Range("C2").Value = "Détails"
$ws.Range("D2").Value = "Framework"
$ws.Range("E2").Value = "Durée"

# ---- Login / Laravel rows ----
$ws.Range("B3").Value = "Login"
$ws.Range("C3").Value = "Creation controller AuthController avec la fonction login ()"
$ws.Range("D3").Value = "Laravel"
$ws.Range("E3").Value = 5

$ws.Range("C4").Value = "Création route /login"
$ws.Range("D4").Value = "Laravel"
$ws.Range("E4").Value = 2

$ws.Range("C5").Value = "Installer la bibliothèque Sanctum pour créer un token"
$ws.Range("D5").Value = "Laravel"
$ws.Range("E5").Value = 3

# ---- Spring Boot rows ----
$ws.Range("C6").Value = "Ajout de dépendance webflux dans pom.xml"
$ws.Range("D6").Value = "Spring Boot"
$ws.Range("E6").Value = 1

$ws.Range("C7").Value = "Création de LaravelAuthService.java pour relier avec laravel"
$ws.Range("D7").Value = "Spring Boot"
$ws.Range("E7").Value = 5

$ws.Range("C8").Value = "Configurer l'url de l'api dans application-properties localhost:80000/api"
$ws.Range("D8").Value = "Spring Boot"
$ws.Range("E8").Value = 1

$ws.Range("C9").Value = "Création de index.html dans lequel se trouve le login en utilisant Thymeleaf"
$ws.Range("D9").Value = "Spring Boot"
$ws.Range("E9").Value = 3

# ---- Debug rows ----
$ws.Range("C10").Value = "Debug"
$ws.Range("D10").Value = "Laravel"
$ws.Range("E10").Value = 60

$ws.Range("C11").Value = "Debug"
$ws.Range("D11").Value = "Spring Boot"
$ws.Range("E11").Value = 120

# ---- Header formatting: bold 14pt font, yellow fill, centered ----
$headerRange = $ws.Range("B2:E2")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 14
$headerRange.Interior.Color = 65535
$headerRange.HorizontalAlignment = -4108

$ws.Rows.Item(2).RowHeight = 18

# ---- Column widths ----
$ws.Columns.Item(2).ColumnWidth = 7.333333333333333
$ws.Columns.Item(3).ColumnWidth = 61.83333333333333
$ws.Columns.Item(4).ColumnWidth = 48.166666666666664

# ---- Page setup ----
$ws.PageSetup.Orientation = 1

# ---- Selection ----
$ws.Range("G10").Select() | Out-Null
